$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '29.401.34'
Set-TextValue "E2" '  +0.05%  '
Set-TextValue "D3" '1.847.57'
Set-TextValue "E4" '  +0.07%  '
Set-TextValue "D5" '240.24'
Set-TextValue "E5" '  -0.17%  '
Set-TextValue "D6" '0.6288'
Set-TextValue "E6" '  -1.54%  '
Set-TextValue "E7" '  +0.03%  '
Set-TextValue "D8" '0.07613'
Set-TextValue "D9" '0.2931'
Set-TextValue "E9" '  -1.20%  '
Set-TextValue "D10" '24.49'
Set-TextValue "E10" '  -1.09%  '
Set-TextValue "D11" '0.07742'
Set-TextValue "E11" '  -0.01%  '
Set-TextValue "D12" '1.848.93'
Set-TextValue "E12" '  -6.86%  '
Set-TextValue "E13" '  +0.17%  '
Set-TextValue "D14" '0.00001089'
Set-TextValue "E14" '  +9.40%  '
Set-TextValue "D15" '0.6795'
Set-TextValue "E15" '  -0.77%  '
Set-TextValue "D16" '83.79'
Set-TextValue "E16" '  +0.65%  '
Set-TextValue "D17" '2.097.55'
Set-TextValue "E17" '  -7.37%  '
Set-TextValue "D18" '6.184'
Set-TextValue "E18" '  +0.02%  '
Set-TextValue "D19" '29.418.61'
Set-TextValue "E19" '  +0.02%  '
Set-TextValue "D20" '228.77'
Set-TextValue "E20" '  -0.42%  '
Set-TextValue "E21" '  -0.21%  '
Set-TextValue "E22" '  +0.04%  '
Set-TextValue "D23" '7.471'
Set-TextValue "E23" '  -1.39%  '
Set-TextValue "E24" '  +0.04%  '
Set-TextValue "D25" '157.21'
Set-TextValue "E25" '  +0.44%  '
Set-TextValue "D26" '0.1397'
Set-TextValue "E26" '  -1.04%  '
Set-TextValue "D27" '8.361'
Set-TextValue "E27" '  -0.40%  '
Set-TextValue "E28" '  -0.25%  '
Set-TextValue "D29" '1.466'
Set-TextValue "E29" '  -0.24%  '
Set-TextValue "D30" '1.302'
Set-TextValue "E31" '  -2.34%  '
Set-TextValue "E32" '  -0.84%  '
Set-TextValue "D33" '4.031'
Set-TextValue "E33" '  -0.10%  '
Set-TextValue "D34" '1.843'
Set-TextValue "E34" '  -0.28%  '
Set-TextValue "E35" '  -0.19%  '
Set-TextValue "D36" '0.7107'
Set-TextValue "E36" '  -0.91%  '
Set-TextValue "E37" '  -0.22%  '
Set-TextValue "D38" '1.234.08'
Set-TextValue "E38" '  -1.51%  '
Set-TextValue "D39" '2.776'
Set-TextValue "E39" '  -0.50%  '
Set-TextValue "E40" '  -0.73%  '
Set-TextValue "D42" '0.9064'
Set-TextValue "E42" '  -0.19%  '
Set-TextValue "D44" '101.98'
Set-TextValue "E44" '  +0.32%  '
Set-TextValue "D46" '0.00000000121'
Set-TextValue "E46" '  +2.96%  '
Set-TextValue "D47" '7.175'
Set-TextValue "E47" '  +1.54%  '
Set-TextValue "D48" '0.4024'
Set-TextValue "E48" '  -0.17%  '
Set-TextValue "D49" '8.954'
Set-TextValue "E49" '  -2.49%  '
Set-TextValue "E50" '  -1.56%  '
Set-TextValue "D51" '0.1122'
Set-TextValue "E51" '  -0.51%  '
